$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 21:35"

# --- Row 4: Estados Unidos ----------------------------------------------
$ws.Range("B4").Value = 1561766
$ws.Range("C4").Value = 11472
$ws.Range("D4").Value = 361174
$ws.Range("E4").Value = 1107780
$ws.Range("G4").Value = 831
$ws.Range("H4").Value = 92812

# --- Row 10: Francia ------------------------------------------------------
$ws.Range("E10").Value = 90224
$ws.Range("H10").Value = 28022

# --- Row 15: Peru -----------------------------------------------------
$ws.Range("B15").Value = 99483
$ws.Range("C15").Value = 4550
$ws.Range("D15").Value = 36524
$ws.Range("E15").Value = 60045
$ws.Range("G15").Value = 125
$ws.Range("H15").Value = 2914

# --- Row 113: Costa Rica -----------------------------------------------
$ws.Range("B113").Value = 882
$ws.Range("C113").Value = 16
$ws.Range("D113").Value = 577
$ws.Range("E113").Value = 295

# --- Rows 139-142: Togo is now inserted ahead of Cabo Verde in the
# sorted-by-total-cases list, pushing Cabo Verde / Isla de Man / Mauricio
# down one row each (Benin above at row 138 and Madagascar below at row
# 143 keep their places and values).
$ws.Range("A139").Value = "Togo"
$ws.Range("B139").Value = 338
$ws.Range("C139").Value = 8
$ws.Range("D139").Value = 107
$ws.Range("E139").Value = 219
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 12

$ws.Range("A140").Value = "Cabo Verde"
$ws.Range("B140").Value = 335
$ws.Range("C140").Value = 7
$ws.Range("D140").Value = 85
$ws.Range("E140").Value = 247
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 3

$ws.Range("A141").Value = "Isla de Man"
$ws.Range("B141").Value = 335
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 300
$ws.Range("E141").Value = 11
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 24

$ws.Range("A142").Value = "Mauricio"
$ws.Range("B142").Value = 332
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 322
$ws.Range("E142").Value = 0
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 10

# --- Row 178: Angola -----------------------------------------------------
$ws.Range("B178").Value = 52
$ws.Range("C178").Value = 2
$ws.Range("D178").Value = 32
$ws.Range("E178").Value = 32
